$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate metal names from French to English.
# Row 3 keeps its own values, just renamed Cuivre -> Copper.
$ws.Range("A3").Value = "Copper"

# Row 4 ("Mercure") becomes "Lead" and takes on the values that were
# in row 5 ("Plomb").
$ws.Range("A4").Value = "Lead"
$ws.Range("B4").Value = 1.24325
$ws.Range("C4").Value = 1.384
$ws.Range("D4").Value = 0.8983020231213873
$ws.Range("E4").Value = 1.11321134124271

# Row 5 ("Plomb") becomes "Mercury" and takes on the values that were
# in row 4 ("Mercure").
$ws.Range("A5").Value = "Mercury"
$ws.Range("B5").Value = 0.1395
$ws.Range("C5").Value = 0.22
$ws.Range("D5").Value = 0.6340909090909091
$ws.Range("E5").Value = 1.577060931899641
